$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1368.4286
$ws.Range("I125").Value = 890
$ws.Range("K125").Value = 8010
$ws.Range("M125").Value = -5550
$ws.Range("H132").Value = 3324512
$ws.Range("I132").Value = 3665046.5
$ws.Range("K132").Value = 10995139.5
$ws.Range("M132").Value = -10992609.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16411.945
$ws.Range("I32").Value = 17609.385
$ws.Range("K32").Value = 17609.385
$ws.Range("M32").Value = -17322.385
$ws.Range("H74").Value = 513.3570999999999
$ws.Range("I74").Value = 596.625
$ws.Range("J74").Value = 402.33334
$ws.Range("K74").Value = 596.625
$ws.Range("L74").Value = 402.33334
$ws.Range("M74").Value = 277.375
$ws.Range("N74").Value = -2150.33334
$ws.Range("H77").Value = 513.3570999999999
$ws.Range("I77").Value = 596.625
$ws.Range("J77").Value = 402.33334
$ws.Range("K77").Value = 2983.125
$ws.Range("L77").Value = 2011.6667
$ws.Range("M77").Value = 1384.875
$ws.Range("N77").Value = -10747.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2090.138
$ws.Range("I105").Value = 1993.5
$ws.Range("J105").Value = 2304.889
$ws.Range("K105").Value = 1993.5
$ws.Range("L105").Value = 2304.889
$ws.Range("M105").Value = -246.5
$ws.Range("N105").Value = -5798.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 444.44446
$ws.Range("I22").Value = 375
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 375
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -25
$ws.Range("N22").Value = -1700
$ws.Range("H132").Value = 3160.2856
$ws.Range("I132").Value = 1979
$ws.Range("J132").Value = 4735.3335
$ws.Range("K132").Value = 5937
$ws.Range("L132").Value = 14206.0005
$ws.Range("M132").Value = -3407
$ws.Range("N132").Value = -19266.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 30000
$ws.Range("J37").Value = 30000
$ws.Range("L37").Value = 90000
$ws.Range("N37").Value = -90224
$ws.Range("H64").Value = 1280
$ws.Range("I64").Value = 864
$ws.Range("J64").Value = 1800
$ws.Range("K64").Value = 2592
$ws.Range("L64").Value = 5400
$ws.Range("M64").Value = -2322
$ws.Range("N64").Value = -5940
$ws.Range("H67").Value = 1280
$ws.Range("I67").Value = 864
$ws.Range("J67").Value = 1800
$ws.Range("K67").Value = 2592
$ws.Range("L67").Value = 5400
$ws.Range("M67").Value = -1656
$ws.Range("N67").Value = -7272
$ws.Range("H75").Value = 5000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 5000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 15000
$ws.Range("M75").ClearContents() | Out-Null
$ws.Range("N75").Value = -16996
$ws.Range("H78").Value = 5000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 5000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 45000
$ws.Range("M78").ClearContents() | Out-Null
$ws.Range("N78").Value = -54984
$ws.Range("H103").Value = 4366.2856
$ws.Range("I103").Value = 250
$ws.Range("J103").Value = 5488.909
$ws.Range("K103").Value = 750
$ws.Range("L103").Value = 16466.727
$ws.Range("M103").Value = 129
$ws.Range("N103").Value = -18224.727
$ws.Range("H112").Value = 1013.5
$ws.Range("I112").Value = 1013.5
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 3040.5
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -1932.5
$ws.Range("N112").ClearContents() | Out-Null
$ws.Range("H114").Value = 2420.375
$ws.Range("I114").Value = 1425.2222
$ws.Range("J114").Value = 3699.8572
$ws.Range("K114").Value = 4275.6666
$ws.Range("L114").Value = 11099.5716
$ws.Range("M114").Value = -1021.6666
$ws.Range("N114").Value = -17607.5716
$ws.Range("H117").Value = 1717.5
$ws.Range("I117").Value = 500
$ws.Range("J117").Value = 1798.6666
$ws.Range("K117").Value = 1500
$ws.Range("L117").Value = 5395.9998
$ws.Range("M117").Value = 1942
$ws.Range("N117").Value = -12279.9998
$ws.Range("H121").Value = 1191160.6
$ws.Range("I121").Value = 390
$ws.Range("J121").Value = 1429314.8
$ws.Range("K121").Value = 1170
$ws.Range("L121").Value = 4287944.4
$ws.Range("M121").Value = 140
$ws.Range("N121").Value = -4290564.4
$ws.Range("H131").Value = 3635.257
$ws.Range("J131").Value = 870.8276
$ws.Range("L131").Value = 2612.4828
$ws.Range("N131").Value = -12692.4828
$ws.Range("H137").Value = 30875172
$ws.Range("I137").Value = 20834020
$ws.Range("J137").Value = 45480490
$ws.Range("K137").Value = 62502060
$ws.Range("L137").Value = 136441470
$ws.Range("M137").Value = -62496960
$ws.Range("N137").Value = -136451670
$ws.Range("H141").Value = 2209.1667
$ws.Range("I141").Value = 2051
$ws.Range("K141").Value = 6153
$ws.Range("M141").Value = -973

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8674.799999999999
$ws.Range("I126").Value = 5159
$ws.Range("J126").Value = 12692.857
$ws.Range("K126").Value = 15477
$ws.Range("L126").Value = 38078.571
$ws.Range("M126").Value = -13007
$ws.Range("N126").Value = -43018.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2554.4546
$ws.Range("I40").Value = 1950
$ws.Range("J40").Value = 2899.8572
$ws.Range("K40").Value = 1950
$ws.Range("L40").Value = 2899.8572
$ws.Range("M40").Value = -1814
$ws.Range("N40").Value = -3171.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents() | Out-Null
$ws.Range("H132").Value = 1414.3928
$ws.Range("I132").Value = 931.1539
$ws.Range("J132").Value = 1833.2
$ws.Range("K132").Value = 2793.4617
$ws.Range("L132").Value = 5499.6
$ws.Range("M132").Value = -263.4616999999998
$ws.Range("N132").Value = -10559.6
$ws.Range("H136").Value = 2273.1177
$ws.Range("I136").Value = 3295.3333
$ws.Range("J136").Value = 1123.125
$ws.Range("K136").Value = 9885.999899999999
$ws.Range("L136").Value = 3369.375
$ws.Range("M136").Value = -7335.999899999999
$ws.Range("N136").Value = -8469.375
